# BRE for Country of document not available in ImmiAccount
# - Adds a new "Country_NotFound" message row to the Messages sheet
# - Shifts all subsequent rows on that sheet down by one
# - Leaves the Messages sheet as the active tab/selection

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("Constants")
$wsMessages  = $wb.Worksheets.Item("Messages")

# Update the selection on the Constants sheet (it is no longer the active tab)
$wsConstants.Range("A27").Select() | Out-Null

# Insert a new row 12 on the Messages sheet for the new Country_NotFound message
$wsMessages.Rows.Item(12).Insert()
$wsMessages.Rows.Item(12).RowHeight = 15

$wsMessages.Range("A12").Value = "Country_NotFound"
$wsMessages.Range("B12").Value = "RPA could not complete this Vevo check. 'Country of Document' value does not match any options on the ImmiAccount website - please check that a valid value has been provided."

# Make Messages the active sheet with B12 selected
$wsMessages.Range("B12").Select() | Out-Null
$wsMessages.Activate()
